$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes $TextValue into $CellAddr as a genuine text cell (not
# auto-coerced to a number), by building it through a TEXT-producing
# formula in a scratch cell and pasting only the resulting value across.
function Set-TextValue {
    param($Sheet, $CellAddr, $TextValue)
    $helper = $Sheet.Range("ZZ1")
    $escaped = $TextValue.Replace('"', '""')
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $Sheet.Range($CellAddr).PasteSpecial(-4163)
    $helper.ClearContents()
}

Set-TextValue $ws "D2" '62.837.27'
Set-TextValue $ws "E2" '  +0.02%  '
Set-TextValue $ws "D3" '3.430.50'
Set-TextValue $ws "E3" '  -0.62%  '
Set-TextValue $ws "E4" '  +0.24%  '
Set-TextValue $ws "D5" '577.70'
Set-TextValue $ws "E5" '  -0.92%  '
Set-TextValue $ws "D6" '146.27'
Set-TextValue $ws "E6" '  -0.49%  '
Set-TextValue $ws "D7" '3.429.44'
Set-TextValue $ws "E7" '  -0.64%  '
Set-TextValue $ws "E8" '  -0.07%  '
Set-TextValue $ws "D9" '0.476'
Set-TextValue $ws "E9" '  -0.25%  '
Set-TextValue $ws "D10" '7.76'
Set-TextValue $ws "E10" '  +1.60%  '
Set-TextValue $ws "E11" '  -0.87%  '
Set-TextValue $ws "D12" '0.403'
Set-TextValue $ws "E12" '  +2.96%  '
Set-TextValue $ws "D13" '4.022.81'
Set-TextValue $ws "E13" '  -0.44%  '
Set-TextValue $ws "E14" '  +2.35%  '
Set-TextValue $ws "E15" '  -0.65%  '
Set-TextValue $ws "D16" '3.430.21'
Set-TextValue $ws "E16" '  -0.18%  '
Set-TextValue $ws "E17" '  -1.17%  '
Set-TextValue $ws "D18" '62.856.22'
Set-TextValue $ws "E18" '  +0.43%  '
Set-TextValue $ws "D19" '6.37'
Set-TextValue $ws "E19" '  +2.15%  '
Set-TextValue $ws "D20" '14.35'
Set-TextValue $ws "E20" '  +0.83%  '
Set-TextValue $ws "D21" '9.18'
Set-TextValue $ws "E21" '  -1.16%  '
Set-TextValue $ws "D22" '383.88'
Set-TextValue $ws "E22" '  -3.17%  '
Set-TextValue $ws "D23" '0.559'
Set-TextValue $ws "E23" '  -0.46%  '
Set-TextValue $ws "D24" '74.22'
Set-TextValue $ws "E24" '  -1.91%  '
Set-TextValue $ws "E25" '  -0.40%  '
Set-TextValue $ws "D26" '3.597.67'
Set-TextValue $ws "E26" '  +0.79%  '
Set-TextValue $ws "D27" '0.0000114'
Set-TextValue $ws "E27" '  -2.53%  '
Set-TextValue $ws "E28" '  -5.06%  '
Set-TextValue $ws "D29" '7.58'
Set-TextValue $ws "E29" '  -0.61%  '
Set-TextValue $ws "E30" '  +0.10%  '
Set-TextValue $ws "D31" '8.09'
Set-TextValue $ws "E31" '  -0.54%  '
Set-TextValue $ws "E32" '  -2.14%  '
Set-TextValue $ws "D33" '0.999'
Set-TextValue $ws "E33" '  -0.08%  '
Set-TextValue $ws "D34" '23.22'
Set-TextValue $ws "E34" '  -2.54%  '
Set-TextValue $ws "D35" '1.32'
Set-TextValue $ws "E35" '  -6.10%  '
Set-TextValue $ws "D36" '5.28'
Set-TextValue $ws "E36" '  -1.22%  '
Set-TextValue $ws "D37" '7.07'
Set-TextValue $ws "E37" '  +0.23%  '
Set-TextValue $ws "E38" '  -0.52%  '
Set-TextValue $ws "D39" '31.58'
Set-TextValue $ws "E39" '  +6.33%  '
Set-TextValue $ws "D40" '168.62'
Set-TextValue $ws "E40" '  -0.22%  '
Set-TextValue $ws "D41" '3.467.14'
Set-TextValue $ws "E41" '  -0.57%  '
Set-TextValue $ws "D42" '0.0765'
Set-TextValue $ws "E42" '  -0.08%  '
Set-TextValue $ws "D43" '0.786'
Set-TextValue $ws "E43" '  -0.72%  '
Set-TextValue $ws "D44" '42.31'
Set-TextValue $ws "E44" '  -1.53%  '
Set-TextValue $ws "D45" '1.71'
Set-TextValue $ws "E45" '  -0.02%  '
Set-TextValue $ws "D46" '1.19'
Set-TextValue $ws "E46" '  +0.24%  '
Set-TextValue $ws "E47" '  -3.27%  '
Set-TextValue $ws "D48" '2.582.68'
Set-TextValue $ws "E48" '  +2.37%  '
Set-TextValue $ws "D49" '2.27'
Set-TextValue $ws "E49" '  +6.16%  '
Set-TextValue $ws "E50" '  +0.75%  '
Set-TextValue $ws "D51" '22.73'
Set-TextValue $ws "E51" '  -1.98%  '

$excel.CutCopyMode = $false
